# Auto-generated COM-interop script applying the newTest2.xlsx "bugs" sheet edit:
#  1) columns N/O/P are rotated (N<-O, O<-P, P<-N) for the header row and rows 2-7
#     so the header reads IssueLink / MergeLink / ContainsTheWordFix, and each row's
#     ContainsTheWordFix numeric flag moves from N to P while IssueLink/MergeLink move
#     into N/O.
#  2) five new data rows (8-12) are appended with full commit/bug metadata.
#  3) the sheet dimension / used range grows from A1:S7 to A1:S12 automatically as the
#     new cells are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rotate columns N/O/P for the header and existing rows 1-7 ---
# new N = old O ("IssueLink" col), new O = old P ("MergeLink" col), new P = old N ("ContainsTheWordFix" col)
$ws.Cells.Item(1,14).Value = "IssueLink"  # N1
$ws.Cells.Item(1,15).Value = "MergeLink"  # O1
$ws.Cells.Item(1,16).Value = "ContainsTheWordFix"  # P1

$ws.Cells.Item(2,14).Value = "https://issues.jenkins-ci.org/browse/JENKINS-42717"  # N2
$ws.Cells.Item(2,15).Value = "https://github.com/jenkinsci/jenkins/pull/2815"  # O2
$ws.Cells.Item(2,16).Value = 2  # P2

$ws.Cells.Item(3,14).Value = "https://issues.jenkins-ci.org/browse/JENKINS-42717"  # N3
$ws.Cells.Item(3,15).Value = "https://github.com/jenkinsci/jenkins/pull/2815"  # O3
$ws.Cells.Item(3,16).Value = 2  # P3

$ws.Cells.Item(4,14).Value = "https://issues.jenkins-ci.org/browse/JENKINS-37599"  # N4
$ws.Cells.Item(4,15).Value = "https://github.com/jenkinsci/jenkins/pull/3583"  # O4
$ws.Cells.Item(4,16).Value = 1  # P4

$ws.Cells.Item(5,14).Value = "https://issues.jenkins-ci.org/browse/JENKINS-37599"  # N5
$ws.Cells.Item(5,15).Value = "https://github.com/jenkinsci/jenkins/pull/3583"  # O5
$ws.Cells.Item(5,16).Value = 1  # P5

$ws.Cells.Item(6,14).Value = ""  # N6 -> empty
$ws.Cells.Item(6,15).Value = "https://github.com/jenkinsci/jenkins/pull/3991"  # O6
$ws.Cells.Item(6,16).Value = 2  # P6

$ws.Cells.Item(7,14).Value = ""  # N7 -> empty
$ws.Cells.Item(7,15).Value = "https://github.com/jenkinsci/jenkins/pull/4556"  # O7
$ws.Cells.Item(7,16).Value = 1  # P7

# --- Step 2: append new rows 8-12 ---

# Row 8
$ws.Cells.Item(8,1).Value = 25385  # A8
$ws.Cells.Item(8,2).Value = "b831acd9854b525d680ca72fd218c848121b9d3f"  # B8
$c8 = @'
[JENKINS-42645] Case insensitive search by default for new and anonymous users (#2801)



* [JENKINS-42645] Case insensitive search by default





* [JENKINS-42960] Search in FixedSet more locale friendly





String.equalsIgnoreCase is safer than toLowerCase when non English


locales are used.





* [JENKINS-42645] Review remarks
'@
$ws.Cells.Item(8,3).Value = $c8  # C8
$ws.Cells.Item(8,4).Value = "test/src/test/java/jenkins/widgets/HistoryPageFilterCaseSensitiveSearchTest.java"  # D8
$ws.Cells.Item(8,5).Value = 4690  # E8
$ws.Cells.Item(8,6).Value = 48  # F8
$ws.Cells.Item(8,7).Value = 35  # G8
$ws.Cells.Item(8,8).Value = 13  # H8
$ws.Range("I2").Copy($ws.Cells.Item(8,9))
$ws.Cells.Item(8,9).Value = 42832.44458332176  # I8
$ws.Cells.Item(8,10).Value = "https://github.com/jenkinsci/jenkins/commit/b831acd9854b525d680ca72fd218c848121b9d3f"  # J8
$ws.Cells.Item(8,11).Value = "https://github.com/jenkinsci/jenkins/raw/b831acd9854b525d680ca72fd218c848121b9d3f/test/src/test/java/jenkins/widgets/HistoryPageFilterCaseSensitiveSearchTest.java"  # K8
$ws.Cells.Item(8,12).Value = 136  # L8
$ws.Cells.Item(8,13).Value = "42645,2801,42960"  # M8
$ws.Cells.Item(8,14).Value = "https://issues.jenkins-ci.org/browse/JENKINS-42645,https://issues.jenkins-ci.org/browse/JENKINS-42960"  # N8
$ws.Cells.Item(8,15).Value = "https://github.com/jenkinsci/jenkins/pull/2801"  # O8
$ws.Cells.Item(8,16).Value = 0  # P8
$ws.Cells.Item(8,17).Value = "improvement,bug"  # Q8
$ws.Cells.Item(8,18).Value = "resolved,resolved"  # R8
$ws.Cells.Item(8,19).Value = "done,fixed"  # S8

# Row 9
$ws.Cells.Item(9,1).Value = 6297  # A9
$ws.Cells.Item(9,2).Value = "312fcd1b9ebb5f4ce396c2b7cc93659edd6301c1"  # B9
$c9 = @'
[JENKINS-54854] Added a warning when cron trigger spent more than a threshold (30s) in its execution (#3802)



* [JENKINS-54854] Added a warning when cron trigger spent more than a threshold (30s) in its execution





* Added a test





* Fix based on feedback





* Added Admin monitor





* Cleanup





* Polishing





* Fix according the feedback provided by Oliver





* Added missing renamed admin monitor class





* Better Web layout





* Limit stacked messages to 10 by default





* Improved replacement





* Make SlowTriggerAdminMonitor#errors thread-safe





* Fixed as requested





* Admin monitor title changed
'@
$ws.Cells.Item(9,3).Value = $c9  # C9
$ws.Cells.Item(9,4).Value = "test/src/test/java/jenkins/triggers/TriggerTest.java"  # D9
$ws.Cells.Item(9,5).Value = 1160  # E9
$ws.Cells.Item(9,6).Value = 118  # F9
$ws.Cells.Item(9,7).Value = 118  # G9
$ws.Cells.Item(9,8).Value = 0  # H9
$ws.Range("I2").Copy($ws.Cells.Item(9,9))
$ws.Cells.Item(9,9).Value = 43673.3347337963  # I9
$ws.Cells.Item(9,10).Value = "https://github.com/jenkinsci/jenkins/commit/312fcd1b9ebb5f4ce396c2b7cc93659edd6301c1"  # J9
$ws.Cells.Item(9,11).Value = "https://github.com/jenkinsci/jenkins/raw/312fcd1b9ebb5f4ce396c2b7cc93659edd6301c1/test/src/test/java/jenkins/triggers/TriggerTest.java"  # K9
$ws.Cells.Item(9,12).Value = 117  # L9
$ws.Cells.Item(9,13).Value = "54854,3802"  # M9
$ws.Cells.Item(9,14).Value = "https://issues.jenkins-ci.org/browse/JENKINS-54854"  # N9
$ws.Cells.Item(9,15).Value = "https://github.com/jenkinsci/jenkins/pull/3802"  # O9
$ws.Cells.Item(9,16).Value = 3  # P9
$ws.Cells.Item(9,17).Value = "improvement"  # Q9
$ws.Cells.Item(9,18).Value = "resolved"  # R9
$ws.Cells.Item(9,19).Value = "fixed"  # S9

# Row 10
$ws.Cells.Item(10,1).Value = 20160  # A10
$ws.Cells.Item(10,2).Value = "a79fdaa4b34b8f7fddb39bed3eabf4763940d11b"  # B10
$c10 = @'
Revert "[JENKINS-46911] createProjectFromXML not recognizing unsafe character…" (#3218)



* Revert "[JENKINS-48447] Fixed HTTP 404 error when clicking on newView sidebar link from an other view. (#3178)"





This reverts commit 6df06fc19a4b7ed015ab5213e2dc8d25beb2f607.





* Revert "[JENKINS-46911] createProjectFromXML not recognizing unsafe character… (#3057)"





This reverts commit ac2a1aaf895020bc80fd951ced748820975df6aa.
'@
$ws.Cells.Item(10,3).Value = $c10  # C10
$ws.Cells.Item(10,4).Value = "test/src/test/java/jenkins/triggers/ReverseBuildTriggerTest.java"  # D10
$ws.Cells.Item(10,5).Value = 3552  # E10
$ws.Cells.Item(10,6).Value = 2  # F10
$ws.Cells.Item(10,7).Value = 1  # G10
$ws.Cells.Item(10,8).Value = 1  # H10
$ws.Range("I2").Copy($ws.Cells.Item(10,9))
$ws.Cells.Item(10,9).Value = 43107.76871527778  # I10
$ws.Cells.Item(10,10).Value = "https://github.com/jenkinsci/jenkins/commit/a79fdaa4b34b8f7fddb39bed3eabf4763940d11b"  # J10
$ws.Cells.Item(10,11).Value = "https://github.com/jenkinsci/jenkins/raw/a79fdaa4b34b8f7fddb39bed3eabf4763940d11b/test/src/test/java/jenkins/triggers/ReverseBuildTriggerTest.java"  # K10
$ws.Cells.Item(10,12).Value = 243  # L10
$ws.Cells.Item(10,13).Value = "46911,3218,48447,3178,3057"  # M10
$ws.Cells.Item(10,14).Value = "https://issues.jenkins-ci.org/browse/JENKINS-46911,https://issues.jenkins-ci.org/browse/JENKINS-48447"  # N10
$ws.Cells.Item(10,15).Value = "https://github.com/jenkinsci/jenkins/pull/3218,https://github.com/jenkinsci/jenkins/pull/3178,https://github.com/jenkinsci/jenkins/pull/3057"  # O10
$ws.Cells.Item(10,16).Value = 1  # P10
$ws.Cells.Item(10,17).Value = "bug,bug"  # Q10
$ws.Cells.Item(10,18).Value = "open,resolved"  # R10
$ws.Cells.Item(10,19).Value = "unresolved,fixed"  # S10

# Row 11
$ws.Cells.Item(11,1).Value = 20358  # A11
$ws.Cells.Item(11,2).Value = "ac2a1aaf895020bc80fd951ced748820975df6aa"  # B11
$c11 = @'
[JENKINS-46911] createProjectFromXML not recognizing unsafe character… (#3057)



* [JENKINS-46911] createProjectFromXML not recognizing unsafe character '/'





* Better place for testCreateProjectCheckGoodName()





* Fix failed test





* Make changes suggested on PR review.





* Remove Failure exception, instead throw IOException. Add javadoc





* [JENKINS-46911] - Add TODO according to the comment from @jtnord.
'@
$ws.Cells.Item(11,3).Value = $c11  # C11
$ws.Cells.Item(11,4).Value = "test/src/test/java/jenkins/triggers/ReverseBuildTriggerTest.java"  # D11
$ws.Cells.Item(11,5).Value = 3548  # E11
$ws.Cells.Item(11,6).Value = 2  # F11
$ws.Cells.Item(11,7).Value = 1  # G11
$ws.Cells.Item(11,8).Value = 1  # H11
$ws.Range("I2").Copy($ws.Cells.Item(11,9))
$ws.Cells.Item(11,9).Value = 43107.61053240741  # I11
$ws.Cells.Item(11,10).Value = "https://github.com/jenkinsci/jenkins/commit/ac2a1aaf895020bc80fd951ced748820975df6aa"  # J11
$ws.Cells.Item(11,11).Value = "https://github.com/jenkinsci/jenkins/raw/ac2a1aaf895020bc80fd951ced748820975df6aa/test/src/test/java/jenkins/triggers/ReverseBuildTriggerTest.java"  # K11
$ws.Cells.Item(11,12).Value = 243  # L11
$ws.Cells.Item(11,13).Value = "46911,3057"  # M11
$ws.Cells.Item(11,14).Value = "https://issues.jenkins-ci.org/browse/JENKINS-46911"  # N11
$ws.Cells.Item(11,15).Value = "https://github.com/jenkinsci/jenkins/pull/3057"  # O11
$ws.Cells.Item(11,16).Value = 1  # P11
$ws.Cells.Item(11,17).Value = "bug"  # Q11
$ws.Cells.Item(11,18).Value = "open"  # R11
$ws.Cells.Item(11,19).Value = "unresolved"  # S11

# Row 12
$ws.Cells.Item(12,1).Value = 20927  # A12
$ws.Cells.Item(12,2).Value = "2ae37219fe635d1a93d1bb9a6ad5d79cc4072489"  # B12
$c12 = @'
Merge pull request #3000 from liketic/JENKINS-46161



[Fix JENKINS-46161] Make ReverseBuildTrigger#getUpstreamProjects null…
'@
$ws.Cells.Item(12,3).Value = $c12  # C12
$ws.Cells.Item(12,4).Value = "test/src/test/java/jenkins/triggers/ReverseBuildTriggerTest.java"  # D12
$ws.Cells.Item(12,5).Value = 4146  # E12
$ws.Cells.Item(12,6).Value = 16  # F12
$ws.Cells.Item(12,7).Value = 16  # G12
$ws.Cells.Item(12,8).Value = 0  # H12
$ws.Range("I2").Copy($ws.Cells.Item(12,9))
$ws.Cells.Item(12,9).Value = 42993.77978009259  # I12
$ws.Cells.Item(12,10).Value = "https://github.com/jenkinsci/jenkins/commit/2ae37219fe635d1a93d1bb9a6ad5d79cc4072489"  # J12
$ws.Cells.Item(12,11).Value = "https://github.com/jenkinsci/jenkins/raw/2ae37219fe635d1a93d1bb9a6ad5d79cc4072489/test/src/test/java/jenkins/triggers/ReverseBuildTriggerTest.java"  # K12
$ws.Cells.Item(12,12).Value = 243  # L12
$ws.Cells.Item(12,13).Value = "3000,46161"  # M12
$ws.Cells.Item(12,14).Value = "https://issues.jenkins-ci.org/browse/JENKINS-46161"  # N12
$ws.Cells.Item(12,15).Value = "https://github.com/jenkinsci/jenkins/pull/3000"  # O12
$ws.Cells.Item(12,16).Value = 1  # P12
$ws.Cells.Item(12,17).Value = "improvement"  # Q12
$ws.Cells.Item(12,18).Value = "closed"  # R12
$ws.Cells.Item(12,19).Value = "fixed"  # S12

